# Update the workbook to reflect data refresh through 2022-02-15
# (commit message: "Add data for 2022-02-23")

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sheet (tab) name: "Through 2022-02-14" -> "Through 2022-02-15"
$ws.Name = "Through 2022-02-15"

# Column header text (row 1, column B): through-date label bump
$ws.Range("B1").Value = "February 2022 (through February 15)"

# Englewood (row 2): B2 3 -> 4
$ws.Range("B2").Value = 4

# Austin (row 3): F3 2 -> 3
$ws.Range("F3").Value = 3

# Auburn Gresham (row 7): B7 1 -> 2
$ws.Range("B7").Value = 2

# North Lawndale (row 8): B8 2 -> 3, J8 1 -> 3
$ws.Range("B8").Value = 3
$ws.Range("J8").Value = 3

# United Center (row 9): D9 2 -> 3
$ws.Range("D9").Value = 3

# Garfield Park (row 15): B15 5 -> 6, J15 1 -> 2
$ws.Range("B15").Value = 6
$ws.Range("J15").Value = 2

# Humboldt Park (row 18): L18 2 -> 3
$ws.Range("L18").Value = 3

# Chicago Lawn (row 22): new value H22 = 1
$ws.Range("H22").Value = 1

# River North (row 42): new value N42 = 1
$ws.Range("N42").Value = 1

# Gage Park (row 66): F66 1 -> 2
$ws.Range("F66").Value = 2
